# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-12-10 08:21:54
#
# Applies the updated "Recorded By" attendance lists (re-ordered lists,
# often with the most-recent recorder moved to the front) and the refreshed
# attendance statistics that follow from a newly recorded MICROBIOLOGY
# session (row 13) on the "Session Analysis Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 2 - ANATOMY session 1 : reorder "Recorded By" list
# ---------------------------------------------------------------------
$ws.Range("G2").Value = "System, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg"

# ---------------------------------------------------------------------
# Row 3 - ANATOMY session 2 : reorder "Recorded By" list
# ---------------------------------------------------------------------
$ws.Range("G3").Value = "System, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"

# ---------------------------------------------------------------------
# Row 4 - ANATOMY session 3 : reorder "Recorded By" list
# ---------------------------------------------------------------------
$ws.Range("G4").Value = "asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"

# ---------------------------------------------------------------------
# Row 6 - ANATOMY session 5 : reorder "Recorded By" list
# ---------------------------------------------------------------------
$ws.Range("G6").Value = "majorelle.magdy@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Mohammedeltanany@med.asu.edu.eg, manar.montaser@med.asu.edu.eg"

# ---------------------------------------------------------------------
# Row 7 - BIOCHEMISTRY LAB/CBL session 1 : reorder "Recorded By" list
# ---------------------------------------------------------------------
$ws.Range("G7").Value = "NadaMohamed@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, lamiaa.ossama@med.asu.edu.eg"

# ---------------------------------------------------------------------
# Row 12 - MICROBIOLOGY session 1 : reorder "Recorded By" list
# ---------------------------------------------------------------------
$ws.Range("G12").Value = "dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"

# ---------------------------------------------------------------------
# Row 13 - MICROBIOLOGY session 2 just got recorded.
# Copy the "Recorded" row look (fill/font) from row 12 onto row 13
# (A:I) without touching row 13's own Year/Group/Subject/Session/Date/
# Time values, then fill in the new recorder, student count and status.
# ---------------------------------------------------------------------
$ws.Range("A12:I12").Copy()
$ws.Range("A13:I13").PasteSpecial(-4122)

$ws.Range("G13").Value = "yassmina.fattoh@med.asu.edu.eg"
$ws.Range("H13").Value = "33/251"
$ws.Range("I13").Value = "Recorded"

# ---------------------------------------------------------------------
# Row 15 - PARASITOLOGY session 2 : reorder "Recorded By" list, and the
# HISTOLOGY group statistics block (columns M:S) recomputed after the
# newly-recorded session shifts a "Missing" session into "Recorded".
# ---------------------------------------------------------------------
$ws.Range("G15").Value = "mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"

$ws.Range("O15").Value = 26
$ws.Range("P15").Value = 2

# Percent columns must stay text (e.g. "89.7%"), not auto-converted to a
# numeric percentage by Excel. Force a text format before assigning, then
# restore the original cell look (fill/font/alignment) via a format-only
# paste from a neighbouring cell that already carries that same style.
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "89.7%"
$ws.Range("L8").Copy()
$ws.Range("L9").PasteSpecial(-4122)

$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "27.6%"
$ws.Range("L8").Copy()
$ws.Range("L10").PasteSpecial(-4122)

$ws.Range("R15").NumberFormat = "@"
$ws.Range("R15").Value = "89.7%"
$ws.Range("M15").Copy()
$ws.Range("R15").PasteSpecial(-4122)

$ws.Range("S15").NumberFormat = "@"
$ws.Range("S15").Value = "27.6%"
$ws.Range("M15").Copy()
$ws.Range("S15").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Row 6 / Row 7 "Class Statistics" numbers (Recorded / Missing sessions)
# ---------------------------------------------------------------------
$ws.Range("L6").Value = 26
$ws.Range("L7").Value = 2

# ---------------------------------------------------------------------
# Row 17 - PARASITOLOGY session 5 : reorder "Recorded By" list
# ---------------------------------------------------------------------
$ws.Range("G17").Value = "mohamed.saleem@med.asu.edu.eg, esraa.sami@med.asu.edu.eg"

# ---------------------------------------------------------------------
# Row 20 - PARASITOLOGY SGD/POS session 2 : reorder "Recorded By" list
# ---------------------------------------------------------------------
$ws.Range("G20").Value = "mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"

# ---------------------------------------------------------------------
# Row 24 : reorder "Recorded By" list
# ---------------------------------------------------------------------
$ws.Range("G24").Value = "youstina.gamil@med.asu.edu.eg, Sarah.Mahdy@med.asu.edu.eg"

# ---------------------------------------------------------------------
# Row 27 : reorder "Recorded By" list
# ---------------------------------------------------------------------
$ws.Range("G27").Value = "nourhan.mostafa@med.asu.edu.eg, hana.amr@med.asu.edu.eg"

# ---------------------------------------------------------------------
# Row 30 : reorder "Recorded By" list
# ---------------------------------------------------------------------
$ws.Range("G30").Value = "shorokmohamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
